$d = $word.ActiveDocument

$replacements = @(
    @("33×17=561", "79×80=6320"),
    @("33×61=2013", "52×77=4004"),
    @("99×42=4158", "75×51=3825"),
    @("25×14=350", "49×80=3920"),
    @("27×35=945", "70×53=3710"),
    @("13×88=1144", "46×63=2898"),
    @("79×74=5846", "33×20=660"),
    @("68×36=2448", "37×50=1850"),
    @("49×29=1421", "33×71=2343"),
    @("47×47=2209", "21×17=357"),
    @("72×99=7128", "14×82=1148"),
    @("58×77=4466", "74×11=814"),
    @("96×74=7104", "37×27=999"),
    @("91×67=6097", "50×18=900"),
    @("29×58=1682", "74×94=6956"),
    @("57×28=1596", "33×59=1947"),
    @("62×28=1736", "12×92=1104"),
    @("84×25=2100", "47×27=1269"),
    @("39×90=3510", "18×92=1656"),
    @("97×38=3686", "32×28=896"),
    @("99×17=1683", "61×75=4575"),
    @("82×84=6888", "81×41=3321"),
    @("86×17=1462", "82×57=4674"),
    @("42×96=4032", "79×14=1106"),
    @("19×36=684", "71×42=2982")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Done applying $($replacements.Count) replacements"
